# RBS 7 Setembro - add Ta-Nb-V calibration sheet, drop the scratch "Folha1" sheet.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the old scratch sheet ("Folha1") that only held two helper formulas.
$oldWs = $wb.Worksheets.Item("Folha1")
$oldWs.Delete() | Out-Null

# Add the new "Ta-Nb-V Calib." sheet (ends up after "Runs", becomes the active tab).
$newWs = $wb.Worksheets.Add()
$newWs.Name = "Ta-Nb-V Calib."

# Row/element labels first (so the shared-string table order matches: V, Nb, Ta, ...).
$newWs.Range("B2").Value = "V"
$newWs.Range("B3").Value = "Nb"
$newWs.Range("B4").Value = "Ta"

# Column headers.
$newWs.Range("C1").Value = "170 º"
$newWs.Range("D1").Value = "160 º"
$newWs.Range("E1").Value = "Mean K factor"

$newWs.Range("E2").Select() | Out-Null
